# Apply the "create user" test-data update:
#  1. LoginCredentials: the valid test login is changed from
#     sameer/joshi12345 to admin/administrator.
#  2. A brand new "CreateUser" worksheet is added (at the end of the
#     workbook) holding the data used to drive a "create user" UI test,
#     including two hyperlinked password cells.

$wb = $excel.ActiveWorkbook

# --- 1. Update the LoginCredentials sheet's valid-login row ---------------
$loginSheet = $wb.Worksheets.Item("LoginCredentials")
$loginSheet.Range("A2").Value = "admin"
$loginSheet.Range("B2").Value = "administrator"
[void]$loginSheet.Range("D13").Select()

# --- 2. Add the new CreateUser worksheet at the end ------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$createUser = $wb.Worksheets.Add($null, $lastSheet)
$createUser.Name = "CreateUser"

# Fill the sheet column by column (header, then value) to mirror how the
# data was originally entered.
$createUser.Range("A1").Value = "Location"
$createUser.Range("A2").Value = "Antioch"

$createUser.Range("B1").Value = "FullName"
$createUser.Range("B2").Value = "Bob"

$createUser.Range("C1").Value = "UserName"
$createUser.Range("C2").Value = "Bob123"

$createUser.Range("D1").Value = "AvailableFilter"
$createUser.Range("D2").Value = "QA"

$createUser.Range("E1").Value = "AuthenticationType"
$createUser.Range("E2").Value = "Title21"

$createUser.Range("F1").Value = "NewPassword"
$createUser.Range("F2").Value = "Bob@123456"

$createUser.Range("G1").Value = "ConfirmPassword"
$createUser.Range("G2").Value = "Bob@123456"

# Password cells are hyperlinked (matches the "mailto:" style hyperlinks
# already used on the CreateEmployee sheet)
$createUser.Hyperlinks.Add($createUser.Range("F2"), "mailto:Bob@123456")
$createUser.Hyperlinks.Add($createUser.Range("G2"), "mailto:Bob@123456")

[void]$createUser.Range("I17").Select()
[void]$createUser.Activate()
